$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 1-22: replace the "layer characteristics" sub-table with the
# new "first sublist" sub-table (corrected test data).
$ws.Range("A1").Value = "first sublist"

$ws.Range("A2").Value = "first sublist"
$ws.Range("B2").Value = "first subnode"

$ws.Range("A3").Value = "first sublist"
$ws.Range("B3").Value = "first subnode"
$ws.Range("C3").Value = "one"

$ws.Range("A4").Value = "first sublist"
$ws.Range("B4").Value = "first subnode"
$ws.Range("C4").Value = "two"

$ws.Range("A5").Value = "first sublist"
$ws.Range("B5").Value = "first subnode"
$ws.Range("C5").Value = "three"

$ws.Range("A6").Value = "first sublist"
$ws.Range("B6").Value = "first subnode"
$ws.Range("C6").Value = "three"
$ws.Range("D6").Value = "fine"

$ws.Range("A7").Value = "first sublist"
$ws.Range("B7").Value = "first subnode"
$ws.Range("C7").Value = "three"
$ws.Range("D7").Value = "medium"

$ws.Range("A8").Value = "first sublist"
$ws.Range("B8").Value = "first subnode"
$ws.Range("C8").Value = "three"
$ws.Range("D8").Value = "coarse"

$ws.Range("A9").Value = "first sublist"
$ws.Range("B9").Value = "first subnode"
$ws.Range("C9").Value = "four"

$ws.Range("A10").Value = "first sublist"
$ws.Range("B10").Value = "first subnode"
$ws.Range("C10").Value = "four"
$ws.Range("D10").Value = "fine"

$ws.Range("A11").Value = "first sublist"
$ws.Range("B11").Value = "first subnode"
$ws.Range("C11").Value = "four"
$ws.Range("D11").Value = "medium"

$ws.Range("A12").Value = "first sublist"
$ws.Range("B12").Value = "first subnode"
$ws.Range("C12").Value = "four"
$ws.Range("D12").Value = "coarse"

$ws.Range("A13").Value = "first sublist"
$ws.Range("B13").Value = "first subnode"
$ws.Range("C13").Value = "five"

$ws.Range("A14").Value = "first sublist"
$ws.Range("B14").Value = "first subnode"
$ws.Range("C14").Value = "six"

$ws.Range("A15").Value = "first sublist"
$ws.Range("B15").Value = "second subnode"

$ws.Range("A16").Value = "first sublist"
$ws.Range("B16").Value = "second subnode"
$ws.Range("C16").Value = "one"

$ws.Range("A17").Value = "first sublist"
$ws.Range("B17").Value = "second subnode"
$ws.Range("C17").Value = "two"

$ws.Range("A18").Value = "first sublist"
$ws.Range("B18").Value = "second subnode"
$ws.Range("C18").Value = "three"

$ws.Range("A19").Value = "first sublist"
$ws.Range("B19").Value = "second subnode"
$ws.Range("C19").Value = "four"

$ws.Range("A20").Value = "first sublist"
$ws.Range("B20").Value = "second subnode"
$ws.Range("C20").Value = "five"

$ws.Range("A21").Value = "first sublist"
$ws.Range("B21").Value = "second subnode"
$ws.Range("C21").Value = "six"

$ws.Range("A22").Value = "first sublist"
$ws.Range("B22").Value = "second subnode"
$ws.Range("C22").Value = "seven"

# Rows 27-29: pad "very" with leading/trailing spaces (error-message test
# data for whitespace handling).
$ws.Range("C27").Value = "     very"
$ws.Range("C28").Value = "   very  "
$ws.Range("C29").Value = "                very     "

# Move the active selection to C5 (was D30).
$ws.Range("C5").Select()
